$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Unprotect("D382")

# Insert a new row at 14 for the BL (Blackline Inc) holding; shifts all following rows down by one.
$ws.Rows.Item(14).Insert()

# Re-write holdings rows 2-38 (symbol, name, sector, weight, percent change)
$ws.Cells.Item(2,1).Value2 = 'ABBV'
$ws.Cells.Item(2,2).Value2 = 'Abbvie Inc'
$ws.Cells.Item(2,3).Value2 = 'Health Care'
$ws.Cells.Item(2,4).Value2 = 0.02997949346239954
$ws.Cells.Item(2,5).Value2 = 0.01623406401907723

$ws.Cells.Item(3,1).Value2 = 'ADBE'
$ws.Cells.Item(3,2).Value2 = 'Adobe Inc'
$ws.Cells.Item(3,3).Value2 = 'Information Technology'
$ws.Cells.Item(3,4).Value2 = 0.03029786428368431
$ws.Cells.Item(3,5).Value2 = -0.003636646506291252

$ws.Cells.Item(4,1).Value2 = 'AMD'
$ws.Cells.Item(4,2).Value2 = 'Advanced Micro Devices Inc'
$ws.Cells.Item(4,3).Value2 = 'Information Technology'
$ws.Cells.Item(4,4).Value2 = 0.02958098982489997
$ws.Cells.Item(4,5).Value2 = 0.02951936419830958

$ws.Cells.Item(5,1).Value2 = 'AMZN'
$ws.Cells.Item(5,2).Value2 = 'Amazon.com Inc'
$ws.Cells.Item(5,3).Value2 = 'Consumer Discretionary'
$ws.Cells.Item(5,4).Value2 = 0.06549463257434331
$ws.Cells.Item(5,5).Value2 = 0.008195664364603505

$ws.Cells.Item(6,1).Value2 = 'KMB'
$ws.Cells.Item(6,2).Value2 = 'Kimberly-Clark Corp'
$ws.Cells.Item(6,3).Value2 = 'Consumer Staples'
$ws.Cells.Item(6,4).Value2 = 0.01677957603062487
$ws.Cells.Item(6,5).Value2 = -0.006461127888194418

$ws.Cells.Item(7,1).Value2 = 'QCOM'
$ws.Cells.Item(7,2).Value2 = 'Qualcomm Inc'
$ws.Cells.Item(7,3).Value2 = 'Information Technology'
$ws.Cells.Item(7,4).Value2 = 0.01572017306331455
$ws.Cells.Item(7,5).Value2 = 0.02091454272863547

$ws.Cells.Item(8,1).Value2 = 'AMGN'
$ws.Cells.Item(8,2).Value2 = 'Amgen Inc'
$ws.Cells.Item(8,3).Value2 = 'Health Care'
$ws.Cells.Item(8,4).Value2 = 0.03053767352044476
$ws.Cells.Item(8,5).Value2 = -0.001389210465385537

$ws.Cells.Item(9,1).Value2 = 'BX'
$ws.Cells.Item(9,2).Value2 = 'Blackstone Group Inc'
$ws.Cells.Item(9,3).Value2 = 'Financials'
$ws.Cells.Item(9,4).Value2 = 0.03122822984841345
$ws.Cells.Item(9,5).Value2 = 0.01018867924528299

$ws.Cells.Item(10,1).Value2 = 'NFLX'
$ws.Cells.Item(10,2).Value2 = 'Netflix Inc'
$ws.Cells.Item(10,3).Value2 = 'Communication Services'
$ws.Cells.Item(10,4).Value2 = 0.03238131750526902
$ws.Cells.Item(10,5).Value2 = -0.07400331168004093

$ws.Cells.Item(11,1).Value2 = 'BABA'
$ws.Cells.Item(11,2).Value2 = 'Alibaba Group Holding Ltd'
$ws.Cells.Item(11,3).Value2 = 'Consumer Discretionary'
$ws.Cells.Item(11,4).Value2 = 0.0316045398382849
$ws.Cells.Item(11,5).Value2 = -0.001914042108926317

$ws.Cells.Item(12,1).Value2 = 'ZG'
$ws.Cells.Item(12,2).Value2 = 'Zillow Group Inc'
$ws.Cells.Item(12,3).Value2 = 'Communication Services'
$ws.Cells.Item(12,4).Value2 = 0.01293123681270152
$ws.Cells.Item(12,5).Value2 = 0.04070473876063163

$ws.Cells.Item(13,1).Value2 = 'CIEN'
$ws.Cells.Item(13,2).Value2 = 'Ciena Corp'
$ws.Cells.Item(13,3).Value2 = 'Information Technology'
$ws.Cells.Item(13,4).Value2 = 0.01489488361788666
$ws.Cells.Item(13,5).Value2 = 0.01384530182757993

$ws.Cells.Item(14,1).Value2 = 'BL'
$ws.Cells.Item(14,2).Value2 = 'Blackline Inc'
$ws.Cells.Item(14,3).Value2 = 'Information Technology'
$ws.Cells.Item(14,4).Value2 = 0.01571428094447522
$ws.Cells.Item(14,5).Value2 = 0.02370953630796158

$ws.Cells.Item(15,1).Value2 = 'DFS'
$ws.Cells.Item(15,2).Value2 = 'Discover Financial Services'
$ws.Cells.Item(15,3).Value2 = 'Financials'
$ws.Cells.Item(15,4).Value2 = 0.007604761381953263
$ws.Cells.Item(15,5).Value2 = 0.01652892561983488

$ws.Cells.Item(16,1).Value2 = 'SYF'
$ws.Cells.Item(16,2).Value2 = 'Synchrony Financial'
$ws.Cells.Item(16,3).Value2 = 'Financials'
$ws.Cells.Item(16,4).Value2 = 0.006994534274160681
$ws.Cells.Item(16,5).Value2 = 0.0154157189790245

$ws.Cells.Item(17,1).Value2 = 'CRM'
$ws.Cells.Item(17,2).Value2 = 'Salesforce.Com Inc'
$ws.Cells.Item(17,3).Value2 = 'Information Technology'
$ws.Cells.Item(17,4).Value2 = 0.03134057291428322
$ws.Cells.Item(17,5).Value2 = 0.01140550973855059

$ws.Cells.Item(18,1).Value2 = 'CTXS'
$ws.Cells.Item(18,2).Value2 = 'Citrix Systems Inc'
$ws.Cells.Item(18,3).Value2 = 'Information Technology'
$ws.Cells.Item(18,4).Value2 = 0.02984004664986889
$ws.Cells.Item(18,5).Value2 = 0.02447147408050965

$ws.Cells.Item(19,1).Value2 = 'AKAM'
$ws.Cells.Item(19,2).Value2 = 'Akamai Technologies Inc'
$ws.Cells.Item(19,3).Value2 = 'Information Technology'
$ws.Cells.Item(19,4).Value2 = 0.03104557416439443
$ws.Cells.Item(19,5).Value2 = 0.001613209337635313

$ws.Cells.Item(20,1).Value2 = 'FB'
$ws.Cells.Item(20,2).Value2 = 'Facebook Inc'
$ws.Cells.Item(20,3).Value2 = 'Communication Services'
$ws.Cells.Item(20,4).Value2 = 0.02972082944535324
$ws.Cells.Item(20,5).Value2 = -0.003898893110854007

$ws.Cells.Item(21,1).Value2 = 'GOOG'
$ws.Cells.Item(21,2).Value2 = 'Alphabet Inc'
$ws.Cells.Item(21,3).Value2 = 'Communication Services'
$ws.Cells.Item(21,4).Value2 = 0.04504780177812362
$ws.Cells.Item(21,5).Value2 = -0.0001482366379931266

$ws.Cells.Item(22,1).Value2 = 'GS'
$ws.Cells.Item(22,2).Value2 = 'Goldman Sachs Group Inc'
$ws.Cells.Item(22,3).Value2 = 'Financials'
$ws.Cells.Item(22,4).Value2 = 0.03259127333991024
$ws.Cells.Item(22,5).Value2 = 0.01021453537423156

$ws.Cells.Item(23,1).Value2 = 'HD'
$ws.Cells.Item(23,2).Value2 = 'Home Depot Inc'
$ws.Cells.Item(23,3).Value2 = 'Consumer Discretionary'
$ws.Cells.Item(23,4).Value2 = 0.03181351365311957
$ws.Cells.Item(23,5).Value2 = 0.006142733670823519

$ws.Cells.Item(24,1).Value2 = 'JBHT'
$ws.Cells.Item(24,2).Value2 = 'J B Hunt Transport Services Inc'
$ws.Cells.Item(24,3).Value2 = 'Industrials'
$ws.Cells.Item(24,4).Value2 = 0.02976168146930588
$ws.Cells.Item(24,5).Value2 = -0.009443487557165731

$ws.Cells.Item(25,1).Value2 = 'RCL'
$ws.Cells.Item(25,2).Value2 = 'Royal Caribbean Cruises Ltd'
$ws.Cells.Item(25,3).Value2 = 'Consumer Discretionary'
$ws.Cells.Item(25,4).Value2 = 0.01433022222911818
$ws.Cells.Item(25,5).Value2 = 0.0446527692117924

$ws.Cells.Item(26,1).Value2 = 'WYNN'
$ws.Cells.Item(26,2).Value2 = 'Wynn Resorts Ltd'
$ws.Cells.Item(26,3).Value2 = 'Consumer Discretionary'
$ws.Cells.Item(26,4).Value2 = 0.01433434671230571
$ws.Cells.Item(26,5).Value2 = 0.03576126274251878

$ws.Cells.Item(27,1).Value2 = 'IBM'
$ws.Cells.Item(27,2).Value2 = 'International Business Machines Corp'
$ws.Cells.Item(27,3).Value2 = 'Information Technology'
$ws.Cells.Item(27,4).Value2 = 0.02984868842416656
$ws.Cells.Item(27,5).Value2 = 0.03901273885350331

$ws.Cells.Item(28,1).Value2 = 'MCK'
$ws.Cells.Item(28,2).Value2 = 'Mckesson Corp'
$ws.Cells.Item(28,3).Value2 = 'Health Care'
$ws.Cells.Item(28,4).Value2 = 0.03056517007502827
$ws.Cells.Item(28,5).Value2 = 0.00776229887420965

$ws.Cells.Item(29,1).Value2 = 'MSFT'
$ws.Cells.Item(29,2).Value2 = 'Microsoft Corp'
$ws.Cells.Item(29,3).Value2 = 'Information Technology'
$ws.Cells.Item(29,4).Value2 = 0.03043397222887267
$ws.Cells.Item(29,5).Value2 = 0.00898319522961355

$ws.Cells.Item(30,1).Value2 = 'NKE'
$ws.Cells.Item(30,2).Value2 = 'Nike Inc'
$ws.Cells.Item(30,3).Value2 = 'Consumer Discretionary'
$ws.Cells.Item(30,4).Value2 = 0.02746139827443407
$ws.Cells.Item(30,5).Value2 = 0.02194949256549439

$ws.Cells.Item(31,1).Value2 = 'NVDA'
$ws.Cells.Item(31,2).Value2 = 'NVIDIA Corp'
$ws.Cells.Item(31,3).Value2 = 'Information Technology'
$ws.Cells.Item(31,4).Value2 = 0.03575632317643341
$ws.Cells.Item(31,5).Value2 = 0.01247425228639676

$ws.Cells.Item(32,1).Value2 = 'PEP'
$ws.Cells.Item(32,2).Value2 = 'PepsiCo Inc'
$ws.Cells.Item(32,3).Value2 = 'Consumer Staples'
$ws.Cells.Item(32,4).Value2 = 0.03147982332285256
$ws.Cells.Item(32,5).Value2 = 0.008715942625763473

$ws.Cells.Item(33,1).Value2 = 'TSM'
$ws.Cells.Item(33,2).Value2 = 'Taiwan Semiconductor Manufacturing Co Ltd'
$ws.Cells.Item(33,3).Value2 = 'Information Technology'
$ws.Cells.Item(33,4).Value2 = 0.02934962595847588
$ws.Cells.Item(33,5).Value2 = 0.02131361461504988

$ws.Cells.Item(34,1).Value2 = 'UNH'
$ws.Cells.Item(34,2).Value2 = 'UnitedHealth Group Inc'
$ws.Cells.Item(34,3).Value2 = 'Health Care'
$ws.Cells.Item(34,4).Value2 = 0.03115202511142487
$ws.Cells.Item(34,5).Value2 = 0.005749880210829117

$ws.Cells.Item(35,1).Value2 = 'SHY'
$ws.Cells.Item(35,2).Value2 = 'Ishares Trust Lehman 1 3yr'
$ws.Cells.Item(35,3).ClearContents()
$ws.Cells.Item(35,4).Value2 = 0.03049878553610523
$ws.Cells.Item(35,5).Value2 = 0.000231830300220226

$ws.Cells.Item(36,1).Value2 = 'JPM'
$ws.Cells.Item(36,2).Value2 = 'JPMorgan Chase & Co'
$ws.Cells.Item(36,3).Value2 = 'Financials'
$ws.Cells.Item(36,4).Value2 = 0.0293172193048596
$ws.Cells.Item(36,5).Value2 = 0.00850807262008435

$ws.Cells.Item(37,1).Value2 = 'TGT'
$ws.Cells.Item(37,2).Value2 = 'Target Corp'
$ws.Cells.Item(37,3).Value2 = 'Consumer Discretionary'
$ws.Cells.Item(37,4).Value2 = 0.03256691924870771
$ws.Cells.Item(37,5).Value2 = -0.002412312442707565

$ws.Cells.Item(38,1).Value2 = 'Total'
$ws.Cells.Item(38,2).ClearContents()
$ws.Cells.Item(38,3).ClearContents()
$ws.Cells.Item(38,4).Value2 = 0.9999999999999999
$ws.Cells.Item(38,5).Value2 = 0.007560374086696964

# Update the confidential footer date (2021-04-09 -> 2021-04-21); this text now lives in row 41
$ws.Cells.Item(41,1).Value2 = '***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.
Model holdings provided as of 2021-04-21 for illustrative purposes only and are subject to change.'

# Restore sheet protection with the original password
$ws.Protect("D382")

Write-Output "edit applied"
